$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ===================== Summary sheet (sheet1) =====================
# Row 2: Name
$ws1.Range("B2").Value = "2/10"

# Row 3: Email
$ws1.Range("B3").Value = "1/10"
$ws1.Range("C3").NumberFormat = "@"
$ws1.Range("C3").Value = "10.0%"
$ws1.Range("C2").Copy()
$ws1.Range("C3").PasteSpecial(-4122)

# Row 4: Mobile (style s4 -> s5)
$ws1.Range("C2").Copy()
$ws1.Range("B4").PasteSpecial(-4122)
$ws1.Range("B4").Value = "0/10"
$ws1.Range("C4").NumberFormat = "@"
$ws1.Range("C4").Value = "0.0%"
$ws1.Range("C2").Copy()
$ws1.Range("C4").PasteSpecial(-4122)

# Row 5: Location (style s5 -> s4)
$ws1.Range("B6").Copy()
$ws1.Range("B5").PasteSpecial(-4122)
$ws1.Range("B5").Value = "1/10"
$ws1.Range("C5").NumberFormat = "@"
$ws1.Range("C5").Value = "10.0%"
$ws1.Range("C2").Copy()
$ws1.Range("C5").PasteSpecial(-4122)

# Row 6: Experience
$ws1.Range("B6").Value = "2/10"

# Row 7: Company
$ws1.Range("B7").Value = "7/10"
$ws1.Range("C7").NumberFormat = "@"
$ws1.Range("C7").Value = "70.0%"
$ws1.Range("C2").Copy()
$ws1.Range("C7").PasteSpecial(-4122)

# Row 9: Overall Wrong Cells
$ws1.Range("B9").Value = "13/60"
$ws1.Range("C9").NumberFormat = "@"
$ws1.Range("C9").Value = "21.7%"
$ws1.Range("A9").Copy()
$ws1.Range("C9").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ===================== Report sheet (sheet2) =====================
# Column widths
$ws2.Columns.Item(3).ColumnWidth = 20.166666666666668
$ws2.Columns.Item(4).ColumnWidth = 30.166666666666668
$ws2.Columns.Item(12).ColumnWidth = 32.166666666666664

# Row 4 fixes: E4 and G4
$ws2.Range("E4").Value = "vibhormalik05@gmail.com"
$ws2.Range("E2").Copy()
$ws2.Range("E4").PasteSpecial(-4122)
$ws2.Range("G4").NumberFormat = "@"
$ws2.Range("G4").Value = "8373992286"
$ws2.Range("E2").Copy()
$ws2.Range("G4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New row 7 ---
$ws2.Range("A7").Value = "vishal_kumar.pdf"
$ws2.Range("B7").Value = "Vishal Kumar"
$ws2.Range("C7").Value = "Seeking Opportunity"
$ws2.Range("C2").Copy()
$ws2.Range("C7").PasteSpecial(-4122)
$ws2.Range("D7").Value = "vishaldeep4357@gmail.com"
$ws2.Range("E7").Value = "vishaldeep4357@gmail.com"
$ws2.Range("E2").Copy()
$ws2.Range("E7").PasteSpecial(-4122)
$ws2.Range("F7").Value = 9546299846
$ws2.Range("G7").NumberFormat = "@"
$ws2.Range("G7").Value = "9546299846"
$ws2.Range("E2").Copy()
$ws2.Range("G7").PasteSpecial(-4122)
$ws2.Range("H7").Value = "New Delhi"
$ws2.Range("I7").Value = "new delhi"
$ws2.Range("E2").Copy()
$ws2.Range("I7").PasteSpecial(-4122)
$ws2.Range("J7").Value = 0
$ws2.Range("K7").NumberFormat = "@"
$ws2.Range("K7").Value = "2.0"
$ws2.Range("C2").Copy()
$ws2.Range("K7").PasteSpecial(-4122)
$ws2.Range("L7").Value = "fresher"
$ws2.Range("M7").Value = ""
$ws2.Range("C2").Copy()
$ws2.Range("M7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New row 8 ---
$ws2.Range("A8").Value = "AshishRanjan.pdf"
$ws2.Range("B8").Value = "Ashish Ranjan"
$ws2.Range("C8").Value = "Ashish Ranjan"
$ws2.Range("E2").Copy()
$ws2.Range("C8").PasteSpecial(-4122)
$ws2.Range("D8").Value = "ashishranjan.ar7@gmail.com"
$ws2.Range("E8").Value = "ashishranjan.ar7@gmail.com"
$ws2.Range("E2").Copy()
$ws2.Range("E8").PasteSpecial(-4122)
$ws2.Range("F8").Value = 8404853652
$ws2.Range("G8").NumberFormat = "@"
$ws2.Range("G8").Value = "8404853652"
$ws2.Range("E2").Copy()
$ws2.Range("G8").PasteSpecial(-4122)
$ws2.Range("H8").Value = "Bengaluru"
$ws2.Range("I8").Value = "bengaluru"
$ws2.Range("E2").Copy()
$ws2.Range("I8").PasteSpecial(-4122)
$ws2.Range("J8").Value = 4
$ws2.Range("K8").NumberFormat = "@"
$ws2.Range("K8").Value = "4.1"
$ws2.Range("E2").Copy()
$ws2.Range("K8").PasteSpecial(-4122)
$ws2.Range("L8").Value = "Mindfire Solutions"
$ws2.Range("M8").Value = "mindfire solutions"
$ws2.Range("E2").Copy()
$ws2.Range("M8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New row 9 ---
$ws2.Range("A9").Value = "Thriveni.docx"
$ws2.Range("B9").Value = "Biyyala Thriveni"
$ws2.Range("C9").Value = "Biyyala Thriveni"
$ws2.Range("E2").Copy()
$ws2.Range("C9").PasteSpecial(-4122)
$ws2.Range("D9").Value = "biyyalathriveni77@gmail.com"
$ws2.Range("E9").Value = "biyyalathriveni77@gmail.com"
$ws2.Range("E2").Copy()
$ws2.Range("E9").PasteSpecial(-4122)
$ws2.Range("F9").Value = 7383665921
$ws2.Range("G9").NumberFormat = "@"
$ws2.Range("G9").Value = "7383665921"
$ws2.Range("E2").Copy()
$ws2.Range("G9").PasteSpecial(-4122)
$ws2.Range("H9").Value = "Hyderabad"
$ws2.Range("I9").Value = "hyderabad"
$ws2.Range("E2").Copy()
$ws2.Range("I9").PasteSpecial(-4122)
$ws2.Range("J9").Value = 2
$ws2.Range("K9").NumberFormat = "@"
$ws2.Range("K9").Value = "2.0"
$ws2.Range("E2").Copy()
$ws2.Range("K9").PasteSpecial(-4122)
$ws2.Range("L9").Value = "UI Sottech Pvt Ltd"
$ws2.Range("M9").Value = ""
$ws2.Range("C2").Copy()
$ws2.Range("M9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New row 10 ---
$ws2.Range("A10").Value = "SumitGulliya.pdf"
$ws2.Range("B10").Value = "Sumit Gulliya"
$ws2.Range("C10").Value = "Sumit Gulliya"
$ws2.Range("E2").Copy()
$ws2.Range("C10").PasteSpecial(-4122)
$ws2.Range("D10").Value = "sgulliya@gmail.com"
$ws2.Range("E10").Value = "sgulliya@gmail.com"
$ws2.Range("E2").Copy()
$ws2.Range("E10").PasteSpecial(-4122)
$ws2.Range("F10").Value = 9711672619
$ws2.Range("G10").NumberFormat = "@"
$ws2.Range("G10").Value = "9711672619"
$ws2.Range("E2").Copy()
$ws2.Range("G10").PasteSpecial(-4122)
$ws2.Range("H10").Value = "New Delhi"
$ws2.Range("I10").Value = "delhi"
$ws2.Range("C2").Copy()
$ws2.Range("I10").PasteSpecial(-4122)
$ws2.Range("J10").Value = 11.9
$ws2.Range("K10").NumberFormat = "@"
$ws2.Range("K10").Value = "11.9"
$ws2.Range("E2").Copy()
$ws2.Range("K10").PasteSpecial(-4122)
$ws2.Range("L10").Value = "Bank Of America Continnum India"
$ws2.Range("M10").Value = ""
$ws2.Range("C2").Copy()
$ws2.Range("M10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New row 11 ---
$ws2.Range("A11").Value = "ShinieMehrotra.pdf"
$ws2.Range("B11").Value = "Shinie Mehrotra"
$ws2.Range("C11").Value = "Shinie Mehrotra"
$ws2.Range("E2").Copy()
$ws2.Range("C11").PasteSpecial(-4122)
$ws2.Range("D11").Value = "shinie.mehrotra.vit@gmail.com"
$ws2.Range("E11").Value = "SHINIEMALHOTRA@ICLOUD.COM"
$ws2.Range("C2").Copy()
$ws2.Range("E11").PasteSpecial(-4122)
$ws2.Range("F11").Value = 9538408670
$ws2.Range("G11").NumberFormat = "@"
$ws2.Range("G11").Value = "9538408670"
$ws2.Range("E2").Copy()
$ws2.Range("G11").PasteSpecial(-4122)
$ws2.Range("H11").Value = "Bengaluru"
$ws2.Range("I11").Value = "bengaluru"
$ws2.Range("E2").Copy()
$ws2.Range("I11").PasteSpecial(-4122)
$ws2.Range("J11").Value = 7.7
$ws2.Range("K11").NumberFormat = "@"
$ws2.Range("K11").Value = "7.7"
$ws2.Range("E2").Copy()
$ws2.Range("K11").PasteSpecial(-4122)
$ws2.Range("L11").Value = "Quantiphi Analytics"
$ws2.Range("M11").Value = ""
$ws2.Range("C2").Copy()
$ws2.Range("M11").PasteSpecial(-4122)
$excel.CutCopyMode = 0
